$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the file "8722a735-fc9e-4e87-ac39-ac760696e3a7.md" ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-10-17 14:36:34"

# ---- zh-cn sheet: row 3 is the same file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-10-17 14:36:11"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0414e959d489b7ff2a6b8d35ce251409d4f07a9b/e2e/8722a735-fc9e-4e87-ac39-ac760696e3a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff2fbecd337ca8d6e8da02b21cefc00f1e9a9734/e2e/8722a735-fc9e-4e87-ac39-ac760696e3a7.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- de-de sheet: row 3 is the same file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-10-17 14:36:34"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0414e959d489b7ff2a6b8d35ce251409d4f07a9b/e2e/8722a735-fc9e-4e87-ac39-ac760696e3a7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff2fbecd337ca8d6e8da02b21cefc00f1e9a9734/e2e/8722a735-fc9e-4e87-ac39-ac760696e3a7.md."
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
